# cryptos.xlsx hourly refresh (GitHub Actions, Tue May  7 13:56:31 UTC 2024).
# Updates the Price (D) / Volume(1h) (E) columns for every listed coin.
# The coin in rows 35/36 (Mantle <-> PEPE) and rows 49/50 (USDe <->
# InjectiveProtocol) swapped rank, so those four rows are rewritten in full
# (Coin, Link, Price, Volume(1h)).
#
# Some Price values ("590.41", "5.90", ...) are plain numeric-looking text in
# this sheet (note the preserved trailing zeros and the multi-dot "thousands"
# values like "63.346.87"). Assigning such a string straight to .Value would let
# Excel auto-convert it to a real number and silently drop that formatting, so
# those are entered with a leading apostrophe to force literal text, exactly as
# typing them into Excel by hand would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.346.87'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '3.063.29'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').Value = '''590.41'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = '''153.62'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.536'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').Value = '3.053.03'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').Value = '''0.156'
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('D11').Value = '''5.90'
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = '''0.450'
$ws.Range('E12').Value = '  -2.59%  '
$ws.Range('E13').Value = '  -3.30%  '
$ws.Range('D14').Value = '''36.47'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '3.571.08'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '''7.16'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '63.286.02'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '3.072.66'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('D20').Value = '''482.91'
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').Value = '''14.45'
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').Value = '''0.705'
$ws.Range('E22').Value = '  -4.50%  '
$ws.Range('D23').Value = '''7.52'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').Value = '''2.40'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '''81.79'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').Value = '''12.77'
$ws.Range('E26').Value = '  -4.09%  '
$ws.Range('D27').Value = '''10.59'
$ws.Range('E27').Value = '  +6.55%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').Value = '''7.55'
$ws.Range('E29').Value = '  +1.03%  '
$ws.Range('D30').Value = '''2.68'
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('D31').Value = '''2.21'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').Value = '''27.19'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').Value = '''0.111'
$ws.Range('E34').Value = '  -3.62%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0822'
$ws.Range('E35').Value = '  -3.99%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').Value = '''1.06'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').Value = '''6.04'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('E38').Value = '  -4.70%  '
$ws.Range('D39').Value = '''2.21'
$ws.Range('E39').Value = '  -2.54%  '
$ws.Range('D40').Value = '''9.26'
$ws.Range('E40').Value = '  -1.84%  '
$ws.Range('D41').Value = '''50.58'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').Value = '''439.89'
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('D43').Value = '''0.289'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('D45').Value = '''0.0362'
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('D46').Value = '''39.98'
$ws.Range('E46').Value = '  +3.07%  '
$ws.Range('D47').Value = '2.816.97'
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('D48').Value = '''132.48'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''25.39'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').Value = '''0.999'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = '''2.22'
$ws.Range('E51').Value = '  -3.19%  '
